$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Cells.Item(6, 1).Value = 9881
$ws.Cells.Item(6, 2).Value = 10000
$ws.Cells.Item(6, 3).Value = 81.47
$ws.Cells.Item(6, 4).Value = 80.5
$ws.Cells.Item(6, 5).Value = $false
$ws.Cells.Item(6, 6).Value = -1.19
$ws.Cells.Item(6, 7).Value = 42607.884189814817
$ws.Cells.Item(6, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(6, 8).Value = $false

# Row 7
$ws.Cells.Item(7, 1).Value = 9835.5499999999993
$ws.Cells.Item(7, 2).Value = 9881
$ws.Cells.Item(7, 3).Value = 80.45
$ws.Cells.Item(7, 4).Value = 80.08
$ws.Cells.Item(7, 5).Value = $false
$ws.Cells.Item(7, 6).Value = -0.46
$ws.Cells.Item(7, 7).Value = 42608.616296296299
$ws.Cells.Item(7, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(7, 8).Value = $false
